# Append new vocabulary term rows (184-197) generated from the updated
# Google Sheet export, as per commit "new .ttl from Google sheet has been
# generated".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $A, $B, $D, $F, $G) {
    if ($A) { $ws.Cells.Item($Row, 1).Value = $A }
    if ($B) { $ws.Cells.Item($Row, 2).Value = $B }
    if ($D) { $ws.Cells.Item($Row, 4).Value = $D }
    if ($F) { $ws.Cells.Item($Row, 6).Value = $F }
    if ($G) { $ws.Cells.Item($Row, 7).Value = $G }
}

Set-Row 184 "cl:10164" "kilometre" "km" "cl:10002" $null
Set-Row 185 "cl:10165" "square kilometre" "km2" "cl:10002" $null
Set-Row 186 "cl:10166" "metre" "m" "cl:10002" $null
Set-Row 187 "cl:10167" "percent" "%" "cl:10002" "http://qudt.org/vocab/unit/PERCENT"
Set-Row 188 "cl:10168" "tonne per hectare per year" "t.har.year-1" "cl:10002" $null
Set-Row 189 "cl:10169" "kilogram per hectar" "kg.har-1" "cl:10002" $null
Set-Row 190 "cl:10170" "tonne per year" "t.year-1" "cl:10002" $null
Set-Row 191 "cl:10171" "milligram per kilogram" "mg.kg-1" "cl:10002" $null
Set-Row 192 "cl:10172" "gram per hectare per year" "g.har.year-1" "cl:10002" $null
Set-Row 193 "cl:10173" $null $null "cl:10002" $null
Set-Row 194 "cl:10174" $null $null "cl:10002" $null
Set-Row 195 "cl:10175" $null $null "cl:10002" $null
Set-Row 196 "cl:10176" $null $null "cl:10002" $null
Set-Row 197 "cl:10177" $null $null $null $null
